$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H corresponds to "saptamana 6" (week 6) attendance counts.
# Update / add the week-6 presence values ("cifrele" - digits of a number lesson).
$ws.Range("H7").Value = 1
$ws.Range("H9").Value = 2
$ws.Range("H10").Value = 1
$ws.Range("H12").Value = 2
$ws.Range("H17").Value = 1
$ws.Range("H19").Value = 2
$ws.Range("H21").Value = 2

# Move the active selection to reflect where the user was last working.
$ws.Range("H19").Select()
